$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Row 11 / column B ("Rule" label for the last rule row) changes from the
# text "R40" to the text "1". A leading apostrophe tells Excel to store the
# value as literal text (matching the shared-string <t>1</t> cell in the
# target workbook) instead of coercing it to a number.
$ws.Range("B11").Value = "'1"
